# Editar materias primas de varios productos (reordenar/combinar ingredientes)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "3.0-harinita,4.0-leche,1.0-huevos,"
$ws.Range("C3").Value = "1.0-harinita,2.0-manzana,5.0-huevos,"
$ws.Range("C4").Value = "5.0-harinita,1.0-vainilla,2.0-huevos,"
$ws.Range("C5").Value = "5.0-harinita,5.0-huevos,"
$ws.Range("C6").Value = "4.0-harinita,2.0-limon,1.0-crema,5.0-merengue,5.0-huevos,"

$wb.Save()
